# Master Acquisitions.xlsx update
# - Update the Inventory/Acquisitions table on Sheet1 with component rows (RevB of micro board)
# - Add a ROHS column (K) to the header and the new data rows
# - Adjust sheet view (active cell, scroll position) and window geometry to match author's session

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Workbook window geometry (best effort; engine may not persist these) ---
$win = $excel.ActiveWindow
$win.Left = 240
$win.Top = 195
$win.Width = 20115
$win.Height = 7875

# --- New "ROHS" header cell, matching style of existing header row (Heading 2 look) ---
$ws.Cells.Item(8, 11).Value = "ROHS"
$ws.Cells.Item(8, 11).Style = $ws.Cells.Item(8, 10).Style

# --- New column widths for columns K & L ---
$ws.Columns.Item(11).ColumnWidth = 11.307291666666666
$ws.Columns.Item(12).ColumnWidth = 11.307291666666666

# --- Row 9: Capacitor ~10 4.7u C203 SMD_0603 ---
$ws.Cells.Item(9, 1).Value = "Capacitor"
$ws.Cells.Item(9, 2).Value = "~10"
$ws.Cells.Item(9, 3).Value = "4.7u"
$ws.Cells.Item(9, 4).Value = "C203"
$ws.Cells.Item(9, 5).Value = "SMD_0603"

# --- Row 10: Microcontroller ---
$ws.Cells.Item(10, 1).Value = "Microcontroller"
$ws.Cells.Item(10, 2).Value = "~6"
$ws.Cells.Item(10, 3).Value = "N/A"
$ws.Cells.Item(10, 4).Value = "U101"
$ws.Cells.Item(10, 5).Value = "TQFP ??"
$ws.Cells.Item(10, 6).Value = "Mouser"
$ws.Cells.Item(10, 7).Value = "556-ATMEGA328P-AU"
$ws.Cells.Item(10, 8).Value = "Atmel"
$ws.Cells.Item(10, 9).Value = "ATMEGA328P-AU"
$ws.Cells.Item(10, 10).Value = "3.64/2.74/2.58"
$ws.Cells.Item(10, 11).Value = "compliant"

# --- Row 11: Resonator ---
$ws.Cells.Item(11, 1).Value = "??????Resonator"
$ws.Cells.Item(11, 2).Value = "~4"
$ws.Cells.Item(11, 3).Value = "16M"
$ws.Cells.Item(11, 4).Value = "X101"
$ws.Cells.Item(11, 5).Value = "CSTCE"
$ws.Cells.Item(11, 6).Value = "Mouser"
$ws.Cells.Item(11, 7).Value = "81-CSTCE16M0V53-R0"
$ws.Cells.Item(11, 8).Value = "Murata"
$ws.Cells.Item(11, 9).Value = "CSTCE16M0V53-R0"
$ws.Cells.Item(11, 10).Value = "0.434/0.35/--"
$ws.Cells.Item(11, 11).Value = "compliant"

# --- Row 12: Resistor 10k R104 ---
$ws.Cells.Item(12, 1).Value = "Resistor"
$ws.Cells.Item(12, 2).Value = "~10"
$ws.Cells.Item(12, 3).Value = "10k"
$ws.Cells.Item(12, 4).Value = "R104"
$ws.Cells.Item(12, 5).Value = "SMD_0603"

# --- Row 13: Resistor 180(ohm, numeric) R106,R107 ---
$ws.Cells.Item(13, 1).Value = "Resistor"
$ws.Cells.Item(13, 2).Value = "~14"
$ws.Cells.Item(13, 3).Value = 180
$ws.Cells.Item(13, 4).Value = "R106,R107"
$ws.Cells.Item(13, 5).Value = "SMD_0603"

# --- Row 14: Capacitor 10n C105 ---
$ws.Cells.Item(14, 1).Value = "Capacitor"
$ws.Cells.Item(14, 2).Value = "~10"
$ws.Cells.Item(14, 3).Value = "10n"
$ws.Cells.Item(14, 4).Value = "C105"
$ws.Cells.Item(14, 5).Value = "SMD_0603"

# --- Row 15: Resistor 0(ohm, numeric) R105,R108 ---
$ws.Cells.Item(15, 1).Value = "Resistor"
$ws.Cells.Item(15, 2).Value = "~15"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = "R105,R108"
$ws.Cells.Item(15, 5).Value = "SMD_0603"

# --- Row 16: Regulator ---
$ws.Cells.Item(16, 1).Value = "Regulator"
$ws.Cells.Item(16, 2).Value = "~6"
$ws.Cells.Item(16, 3).Value = "N/A"
$ws.Cells.Item(16, 4).Value = "U103"
$ws.Cells.Item(16, 5).Value = "SOT-23-5"
$ws.Cells.Item(16, 6).Value = "Mouser"
$ws.Cells.Item(16, 7).Value = "595-LP2985-33DBVR"
$ws.Cells.Item(16, 8).Value = "TI"
$ws.Cells.Item(16, 9).Value = "LP2985-33DBVR"
$ws.Cells.Item(16, 10).Value = "0.578/0.43/0.317"
$ws.Cells.Item(16, 11).Value = "compliant"

# --- Row 17: Pin Headers ---
$ws.Cells.Item(17, 1).Value = "Pin Headers"
$ws.Cells.Item(17, 2).Value = "~5"
$ws.Cells.Item(17, 3).Value = "N/A"
$ws.Cells.Item(17, 4).Value = "N/A"
$ws.Cells.Item(17, 5).Value = "N/A"
$ws.Cells.Item(17, 6).Value = "Mouser"
$ws.Cells.Item(17, 7).Value = "517-9611106404AR"
$ws.Cells.Item(17, 8).Value = "3M"
$ws.Cells.Item(17, 9).Value = "961110-6404-AR"
$ws.Cells.Item(17, 10).Value = "'--/0.437/0.412"
$ws.Cells.Item(17, 11).Value = "compliant"

# --- View state: scroll so column H is at the left edge, select J18 ---
$ws.Activate()
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("J18").Select()
